# "updated spreadsheet corrected with NA in comments"
#
# Fill the previously-empty "comments" column (M) cells for the rows that
# didn't have a reviewer comment yet with the literal text "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 6, 7, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29)

foreach ($r in $rows) {
    $ws.Range("M$r").Value = "NA"
}

# Move the selection onto the newly-filled range, matching the author's
# resulting cursor position/selection (M14 active, M14:M29 selected).
$ws.Range("M14:M29").Select()
